$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.028.63"
$ws.Range("E2").Value = "  +1.60%  "

# Row 3
$ws.Range("D3").Value = "3.133.54"
$ws.Range("E3").Value = "  +0.34%  "

# Row 4
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").Value = "'591.90"
$ws.Range("E5").Value = "  +0.46%  "

# Row 6
$ws.Range("D6").Value = "'147.35"
$ws.Range("E6").Value = "  +1.90%  "

# Row 7
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("D8").Value = "3.130.93"
$ws.Range("E8").Value = "  +0.61%  "

# Row 9
$ws.Range("E9").Value = "  -0.20%  "

# Row 10
$ws.Range("D10").Value = "'0.162"
$ws.Range("E10").Value = "  +12.05%  "

# Row 11
$ws.Range("D11").Value = "'5.75"
$ws.Range("E11").Value = "  -0.45%  "

# Row 12
$ws.Range("D12").Value = "'0.469"
$ws.Range("E12").Value = "  -0.71%  "

# Row 13
$ws.Range("E13").Value = "  +3.99%  "

# Row 14
$ws.Range("D14").Value = "'37.50"
$ws.Range("E14").Value = "  +4.78%  "

# Row 15
$ws.Range("E15").Value = "  -1.02%  "

# Row 16
$ws.Range("D16").Value = "3.653.66"
$ws.Range("E16").Value = "  +0.53%  "

# Row 17
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "'7.19"
$ws.Range("E17").Value = "  -2.33%  "

# Row 18
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "63.860.54"
$ws.Range("E18").Value = "  +1.52%  "

# Row 19
$ws.Range("D19").Value = "3.139.74"
$ws.Range("E19").Value = "  +0.98%  "

# Row 20
$ws.Range("D20").Value = "'468.44"
$ws.Range("E20").Value = "  +2.81%  "

# Row 21
$ws.Range("D21").Value = "'14.38"
$ws.Range("E21").Value = "  +1.35%  "

# Row 22
$ws.Range("D22").Value = "'0.736"
$ws.Range("E22").Value = "  -0.38%  "

# Row 23
$ws.Range("D23").Value = "'7.58"

# Row 24
$ws.Range("D24").Value = "'13.33"
$ws.Range("E24").Value = "  -3.91%  "

# Row 25
$ws.Range("D25").Value = "'82.54"
$ws.Range("E25").Value = "  -0.07%  "

# Row 26
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.27%  "

# Row 27
$ws.Range("D27").Value = "'9.02"
$ws.Range("E27").Value = "  +8.36%  "

# Row 28
$ws.Range("E28").Value = "  +0.34%  "

# Row 29
$ws.Range("E29").Value = "  -1.55%  "

# Row 30
$ws.Range("E30").Value = "  -0.04%  "

# Row 31
$ws.Range("D31").Value = "'6.88"

# Row 32
$ws.Range("D32").Value = "'27.23"
$ws.Range("E32").Value = "  -0.19%  "

# Row 33
$ws.Range("E33").Value = "  -5.08%  "

# Row 34
$ws.Range("D34").Value = "0.0₃0882"
$ws.Range("E34").Value = "  +9.74%  "

# Row 35
$ws.Range("D35").Value = "'2.40"
$ws.Range("E35").Value = "  +8.08%  "

# Row 36
$ws.Range("E36").Value = "  +0.97%  "

# Row 37
$ws.Range("E37").Value = "  +12.40%  "

# Row 38
$ws.Range("E38").Value = "  +0.36%  "

# Row 39
$ws.Range("D39").Value = "'51.04"
$ws.Range("E39").Value = "  +0.65%  "

# Row 40
$ws.Range("D40").Value = "'458.03"
$ws.Range("E40").Value = "  +6.93%  "

# Row 41
$ws.Range("D41").Value = "'8.74"
$ws.Range("E41").Value = "  -1.83%  "

# Row 42
$ws.Range("D42").Value = "'0.0373"
$ws.Range("E42").Value = "  +0.17%  "

# Row 43
$ws.Range("D43").Value = "2.905.35"
$ws.Range("E43").Value = "  -1.22%  "

# Row 44
$ws.Range("D44").Value = "'0.278"
$ws.Range("E44").Value = "  -0.33%  "

# Row 45
$ws.Range("E45").Value = "  +1.09%  "

# Row 46
$ws.Range("E46").Value = "  -0.44%  "

# Row 47
$ws.Range("D47").Value = "'126.48"
$ws.Range("E47").Value = "  +1.88%  "

# Row 48
$ws.Range("D48").Value = "'35.70"
$ws.Range("E48").Value = "  +1.28%  "

# Row 50
$ws.Range("E50").Value = "  -0.50%  "

# Row 51
$ws.Range("D51").Value = "'24.83"
$ws.Range("E51").Value = "  -0.44%  "
